$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 46: date 46029 (2026-01-07), value 35
$ws.Range("A46").Value = 46029
$ws.Range("B46").Value = 35

# Add new row 47: date 46028 (2026-01-06), value 44
$ws.Range("A47").Value = 46028
$ws.Range("B47").Value = 44

# Match formatting of the preceding date cell (re-use existing date style)
$ws.Range("A45").Copy()
$ws.Range("A46:A47").PasteSpecial(-4122)

# Update the selection to match the new active cell/selection state
$ws.Range("G48").Select()
